$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old values while keeping existing formatting on A1:E1 (header) and A2:E11 (data)
$ws.Range("A1:E11").ClearContents()

# Extend the header formatting (bold, border, centered) from A1 across the new header columns F1:P1
$ws.Range("A1").Copy()
$ws.Range("F1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @("One Room RANSAC", "Time", "One Room Refined", "Time", "Two Rooms RANSAC", "Time", "Two Rooms Refined", "Time", "Three Rooms RANSAC", "Time", "Three Rooms Refined", "Time", "Four Rooms RANSAC", "Time", "Four Rooms Refined", "Time")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

$data = @{}
$data[2] = @(0.2518898123372789, 0.1000039577484131, 0.1482372422167708, 0.05199837684631348, 0.2261681963224231, 0.0800025463104248, 0.1633780183892594, 0.02500176429748535, 0.2436729525948772, 0.1100020408630371, 0.1473568526337741, 0.1060006618499756, 0.2061270982166633, 0.09900116920471191, 0.1458469042296603, 0.09400224685668945)
$data[3] = @(0.2764373150518605, 0.1010017395019531, 0.1472837945855425, 0.02199602127075195, 0.2604920049314019, 0.09000301361083984, 0.1633779939964034, 0.03400063514709473, 0.3047929911143028, 0.1120002269744873, 0.1473565875268908, 0.1149978637695312, 0.2250521677412161, 0.104001522064209, 0.1458469768804727, 0.1719973087310791)
$data[4] = @(0.3498953023264292, 0.1089985370635986, 0.1589493327675606, 0.05099797248840332, 0.2236084222637152, 0.08100247383117676, 0.1629723345984794, 0.03000020980834961, 0.218573030925764, 0.1090028285980225, 0.1474637758030676, 0.1060028076171875, 0.2463742064013847, 0.09799885749816895, 0.1458469768149112, 0.1659984588623047)
$data[5] = @(0.3001472564275787, 0.1019971370697021, 0.1482372422167708, 0.0579979419708252, 0.2624189214276322, 0.08099794387817383, 0.1629725749536803, 0.03200197219848633, 0.2185730309257643, 0.1140010356903076, 0.1474637758030677, 0.1010012626647949, 0.1981651280651061, 0.1120011806488037, 0.1458472395410905, 0.1030004024505615)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
